$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2-62, column G
$kValues = @(4,2,0,1,3,1,0,0,2,1,0,3,0,1,0,0,2,1,2,0,1,1,2,1,1,3,0,3,0,1,3,3,1,2,0,3,1,3,4,2,3,2,1,1,2,2,1,0,2,4,4,2,1,0,1,2,1,0,1,2,2)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
